$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $result = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
    if (-not $result) {
        Write-Output "FAILED: $old -> $new"
    }
}

Replace-Text "2025-02-28 Friday" "2025-03-01 Saturday"
Replace-Text "361×5=1805" "150×7=1050"
Replace-Text "358×9=3222" "715×6=4290"
Replace-Text "382×6=2292" "518×7=3626"
Replace-Text "732×4=2928" "289×2=578"
Replace-Text "577×8=4616" "456×9=4104"
Replace-Text "685×3=2055" "329×7=2303"
Replace-Text "219×7=1533" "326×9=2934"
Replace-Text "558×9=5022" "481×7=3367"
Replace-Text "195×8=1560" "398×9=3582"
Replace-Text "443×3=1329" "725×4=2900"
Replace-Text "737×6=4422" "168×8=1344"
Replace-Text "998×8=7984" "738×2=1476"
Replace-Text "164×2=328" "827×9=7443"
Replace-Text "795×7=5565" "698×6=4188"
Replace-Text "536×4=2144" "287×6=1722"
Replace-Text "673×2=1346" "957×8=7656"
Replace-Text "246×9=2214" "265×9=2385"
Replace-Text "723×9=6507" "898×5=4490"
Replace-Text "112×7=784" "606×8=4848"
Replace-Text "251×5=1255" "911×6=5466"
Replace-Text "149×2=298" "466×6=2796"
Replace-Text "520×9=4680" "321×8=2568"
Replace-Text "746×7=5222" "498×4=1992"
Replace-Text "906×7=6342" "842×7=5894"
Replace-Text "118×2=236" "856×8=6848"

Write-Output "Done"
